# Base donation map setup & control panel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-sort the "The tree has been deleted." (637) row: it currently sits at the
#    bottom of the table (row 67) but numerically belongs right before row 57 (638).
#    Insert a fresh row at 57, populate it, then drop the now-duplicated old row
#    (which has shifted down to row 68).
$ws.Rows.Item(57).Insert()
$ws.Range("A57").Value2 = 637
$ws.Range("B57").Value2 = "Tree"
$ws.Range("C57").Value2 = "The tree has been deleted."
$ws.Range("D57").Value2 = "Response"
$ws.Rows.Item(68).Delete()

# 2) Insert 4 new rows for the "Location" category right after the Tree rows
#    (before the Authentification block) and fill them in, in the same order
#    the author originally typed them.
$ws.Rows.Item(60).Resize(4).Insert()

$ws.Range("A60").Value2 = 664
$ws.Range("B60").Value2 = "Location"
$ws.Range("C60").Value2 = "The location has been updated."
$ws.Range("D60").Value2 = "Response"

$ws.Range("A61").Value2 = 665
$ws.Range("B61").Value2 = "Location"
$ws.Range("C61").Value2 = "New location has been created."
$ws.Range("D61").Value2 = "Response"

$ws.Range("A63").Value2 = 667
$ws.Range("B63").Value2 = "Location"
$ws.Range("C63").Value2 = "The location has been deleted."
$ws.Range("D63").Value2 = "Response"

$ws.Range("A62").Value2 = 666
$ws.Range("B62").Value2 = "Location"
$ws.Range("C62").Value2 = "Press the button to delete this location."
$ws.Range("D62").Value2 = "Message"

# 3) Re-select the newly added block, mirroring the author's on-screen selection
#    (the existing AutoFilter on A2:D2 is left as-is so it stays turned on).
$ws.Range("A60:D63").Select()
